$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values
$ws.Range("D2").Value = "96.264.40"
$ws.Range("E2").Value = "  +4.84%  "
$ws.Range("D3").Value = "3.640.40"
$ws.Range("E3").Value = "  +9.41%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("E5").Value = "  +3.99%  "
$ws.Range("E6").Value = "  +4.18%  "
$ws.Range("E7").Value = "  +6.65%  "
$ws.Range("E8").Value = "  +4.69%  "
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("E10").Value = "  +6.86%  "
$ws.Range("D11").Value = "3.639.24"
$ws.Range("E11").Value = "  +9.43%  "
$ws.Range("E12").Value = "  +2.35%  "
$ws.Range("E14").Value = "  +5.85%  "
$ws.Range("D15").Value = "4.320.62"
$ws.Range("E15").Value = "  +9.43%  "
$ws.Range("D16").Value = "96.217.84"
$ws.Range("E16").Value = "  +4.98%  "
$ws.Range("E17").Value = "  +5.34%  "
$ws.Range("D18").Value = "3.642.63"
$ws.Range("E18").Value = "  +9.85%  "
$ws.Range("E19").Value = "  +23.21%  "
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("E21").Value = "  +6.40%  "
$ws.Range("E22").Value = "  +10.77%  "
$ws.Range("E23").Value = "  +5.17%  "
$ws.Range("E24").Value = "  +0.49%  "
$ws.Range("E25").Value = "  +8.44%  "
$ws.Range("E26").Value = "  +8.00%  "
$ws.Range("E27").Value = "  +8.63%  "
$ws.Range("E28").Value = "  +5.92%  "
$ws.Range("E29").Value = "  +20.64%  "
$ws.Range("E30").Value = "  +5.44%  "
$ws.Range("E31").Value = "  +3.42%  "
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("E33").Value = "  +6.04%  "
$ws.Range("E34").Value = "  +0.70%  "
$ws.Range("E35").Value = "  +10.97%  "
$ws.Range("E36").Value = "  +9.42%  "
$ws.Range("E37").Value = "  +2.96%  "
$ws.Range("E38").Value = "  +7.54%  "
$ws.Range("E39").Value = "  +8.44%  "
$ws.Range("E40").Value = "  +3.36%  "
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("E41").Value = "  +8.63%  "
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("E43").Value = "  +4.65%  "
$ws.Range("E44").Value = "  +4.57%  "
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("E45").Value = "  +6.33%  "
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("E46").Value = "  +0.58%  "
$ws.Range("E47").Value = "  +5.94%  "
$ws.Range("E48").Value = "  -2.29%  "
$ws.Range("E49").Value = "  +4.64%  "
$ws.Range("E50").Value = "  +3.49%  "
$ws.Range("E51").Value = "  +4.12%  "

# Numeric-looking Price values: force text storage (leading apostrophe)
# then reset the cell style back to Normal so no stray quote-prefix style lingers
$c = $ws.Range("D5")
$c.Value = "'240.06"
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.Value = "'640.20"
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.Value = "'1.50"
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.Value = "'0.402"
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.Value = "'43.45"
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.Value = "'6.36"
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.Value = "'13.41"
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.Value = "'8.02"
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.Value = "'18.36"
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.Value = "'0.499"
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.Value = "'517.13"
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.Value = "'97.37"
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.Value = "'12.51"
$c.Style = "Normal"
$c = $ws.Range("D29")
$c.Value = "'3.13"
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.Value = "'11.61"
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.Value = "'1.00"
$c.Style = "Normal"
$c = $ws.Range("D35")
$c.Value = "'31.22"
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.Value = "'0.574"
$c.Style = "Normal"
$c = $ws.Range("D37")
$c.Value = "'572.04"
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.Value = "'7.87"
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.Value = "'0.939"
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.Value = "'1.00"
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.Value = "'5.73"
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.Value = "'23.80"
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.Value = "'54.10"
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.Value = "'8.20"
$c.Style = "Normal"
